$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 110.8125
$ws.Range("I9").Value = 112.61539
$ws.Range("K9").Value = 112.61539
$ws.Range("M9").Value = 56.38461
$ws.Range("H19").Value = 21739930
$ws.Range("I19").Value = 783.0909
$ws.Range("K19").Value = 783.0909
$ws.Range("M19").Value = -608.0909
$ws.Range("H48").Value = 4998
$ws.Range("I48").Value = 4998
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 14994
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -14702
$ws.Range("N48").ClearContents() | Out-Null
$ws.Range("H56").Value = 4998
$ws.Range("I56").Value = 4998
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 14994
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -14460
$ws.Range("N56").ClearContents() | Out-Null
$ws.Range("H57").Value = 50709
$ws.Range("J57").Value = 50709
$ws.Range("L57").Value = 152127
$ws.Range("N57").Value = -153125
$ws.Range("H98").Value = 1344.4
$ws.Range("I98").Value = 1280.4
$ws.Range("J98").Value = 1600.4
$ws.Range("K98").Value = 1280.4
$ws.Range("L98").Value = 1600.4
$ws.Range("M98").Value = 217.5999999999999
$ws.Range("N98").Value = -4596.4
$ws.Range("H122").Value = 1344.4
$ws.Range("I122").Value = 1280.4
$ws.Range("J122").Value = 1600.4
$ws.Range("K122").Value = 3841.2
$ws.Range("L122").Value = 4801.200000000001
$ws.Range("M122").Value = -1391.2
$ws.Range("N122").Value = -9701.200000000001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H40").Value = 55555
$ws.Range("J40").Value = 55555
$ws.Range("L40").Value = 55555
$ws.Range("N40").Value = -55907
$ws.Range("H110").Value = 1635.9584
$ws.Range("I110").Value = 1432.4615
$ws.Range("K110").Value = 1432.4615
$ws.Range("M110").Value = 612.5385000000001
$ws.Range("H122").Value = 2304.6667
$ws.Range("I122").Value = 2349.2727
$ws.Range("K122").Value = 7047.8181
$ws.Range("M122").Value = -4597.8181

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 14986.333
$ws.Range("J88").Value = 14986.333
$ws.Range("L88").Value = 14986.333
$ws.Range("N88").Value = -15798.333
$ws.Range("H91").Value = 14986.333
$ws.Range("J91").Value = 14986.333
$ws.Range("L91").Value = 14986.333
$ws.Range("N91").Value = -17794.333

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2501019.2
$ws.Range("I31").Value = 4167345.8
$ws.Range("J31").Value = 1529.35
$ws.Range("K31").Value = 4167345.8
$ws.Range("L31").Value = 1529.35
$ws.Range("M31").Value = -4167050.8
$ws.Range("N31").Value = -2119.35
$ws.Range("H34").Value = 2501019.2
$ws.Range("I34").Value = 4167345.8
$ws.Range("J34").Value = 1529.35
$ws.Range("K34").Value = 4167345.8
$ws.Range("L34").Value = 1529.35
$ws.Range("M34").Value = -4167143.8
$ws.Range("N34").Value = -1933.35
$ws.Range("H80").Value = 99999
$ws.Range("J80").Value = 99999
$ws.Range("L80").Value = 99999
$ws.Range("N80").Value = -102245
$ws.Range("H83").Value = 99999
$ws.Range("J83").Value = 99999
$ws.Range("L83").Value = 299997
$ws.Range("N83").Value = -311229
$ws.Range("H88").Value = 26114
$ws.Range("J88").Value = 19171
$ws.Range("L88").Value = 19171
$ws.Range("N88").Value = -19983
$ws.Range("H91").Value = 26114
$ws.Range("J91").Value = 19171
$ws.Range("L91").Value = 19171
$ws.Range("N91").Value = -21979
$ws.Range("H122").Value = 1920.75
$ws.Range("I122").Value = 1799.5
$ws.Range("K122").Value = 5398.5
$ws.Range("M122").Value = -2948.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1572.2368
$ws.Range("J5").Value = 2741.9092
$ws.Range("L5").Value = 8225.7276
$ws.Range("N5").Value = -8449.7276
$ws.Range("H107").Value = 960.4545000000001
$ws.Range("I107").Value = 626.5
$ws.Range("J107").Value = 976.3570999999999
$ws.Range("K107").Value = 1879.5
$ws.Range("L107").Value = 2929.0713
$ws.Range("M107").Value = 40.5
$ws.Range("N107").Value = -6769.0713
$ws.Range("H132").Value = 2071
$ws.Range("I132").Value = 1632
$ws.Range("K132").Value = 14688
$ws.Range("M132").Value = -12158
$ws.Range("H135").Value = 1572.2368
$ws.Range("J135").Value = 2741.9092
$ws.Range("L135").Value = 24677.1828
$ws.Range("N135").Value = -29747.1828

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 38114.125
$ws.Range("I46").Value = 13313.667
$ws.Range("K46").Value = 13313.667
$ws.Range("M46").Value = -13157.667
$ws.Range("H53").Value = 9000
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents() | Out-Null
$ws.Range("H74").Value = 90000
$ws.Range("J74").Value = 90000
$ws.Range("L74").Value = 90000
$ws.Range("N74").Value = -91872
$ws.Range("H77").Value = 90000
$ws.Range("J77").Value = 90000
$ws.Range("L77").Value = 270000
$ws.Range("N77").Value = -279360
$ws.Range("H102").Value = 4999.3335
$ws.Range("I102").Value = 4999.5
$ws.Range("J102").Value = 4999
$ws.Range("K102").Value = 4999.5
$ws.Range("L102").Value = 4999
$ws.Range("M102").Value = -3377.5
$ws.Range("N102").Value = -8243
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
$ws.Range("H123").Value = 36326
$ws.Range("J123").Value = 36326
$ws.Range("L123").Value = 36326
$ws.Range("N123").Value = -41226

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2985.7144
$ws.Range("I7").Value = 2985.7144
$ws.Range("K7").Value = 2985.7144
$ws.Range("M7").Value = -2873.7144
$ws.Range("H40").Value = 2038.4
$ws.Range("J40").Value = 1798
$ws.Range("L40").Value = 1798
$ws.Range("N40").Value = -2070
$ws.Range("H46").Value = 1574.3077
$ws.Range("I46").Value = 1150.3334
$ws.Range("J46").Value = 1937.7142
$ws.Range("K46").Value = 1150.3334
$ws.Range("L46").Value = 1937.7142
$ws.Range("M46").Value = -962.3334
$ws.Range("N46").Value = -2313.7142
$ws.Range("H80").Value = 44876.4
$ws.Range("J80").Value = 44876.4
$ws.Range("L80").Value = 44876.4
$ws.Range("N80").Value = -47122.4
$ws.Range("H83").Value = 44876.4
$ws.Range("J83").Value = 44876.4
$ws.Range("L83").Value = 134629.2
$ws.Range("N83").Value = -145861.2
$ws.Range("H126").Value = 2985.7144
$ws.Range("I126").Value = 2985.7144
$ws.Range("K126").Value = 8957.143199999999
$ws.Range("M126").Value = -6487.143199999999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3946
$ws.Range("J62").Value = 3999
$ws.Range("L62").Value = 3999
$ws.Range("N62").Value = -5247
$ws.Range("H65").Value = 3946
$ws.Range("J65").Value = 3999
$ws.Range("L65").Value = 19995
$ws.Range("N65").Value = -26235
$ws.Range("H96").Value = 2956.6
$ws.Range("I96").Value = 2698.25
$ws.Range("J96").Value = 3990
$ws.Range("K96").Value = 2698.25
$ws.Range("L96").Value = 3990
$ws.Range("M96").Value = -1325.25
$ws.Range("N96").Value = -6736
$ws.Range("H126").Value = 5349.8
$ws.Range("I126").Value = 5287.6665
$ws.Range("K126").Value = 15862.9995
$ws.Range("M126").Value = -13392.9995
$ws.Range("H130").Value = 96900
$ws.Range("J130").Value = 96900
$ws.Range("L130").Value = 96900
$ws.Range("N130").Value = -106940
$ws.Range("H136").Value = 2301.6785
$ws.Range("I136").Value = 2025.2291
$ws.Range("K136").Value = 6075.6873
$ws.Range("M136").Value = -3525.6873
